$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain text (they contain numeric-looking
# strings like "1.009" that Excel would otherwise auto-convert to numbers).
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range('D2').Value = '26.235.52'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '1.659.40'
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('D4').Value = '1.009'
$ws.Range('E4').Value = '  +0.71%  '
$ws.Range('D5').Value = '218.09'
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('D6').Value = '0.5294'
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D7').Value = '1.009'
$ws.Range('E7').Value = '  +0.63%  '
$ws.Range('D8').Value = '0.2631'
$ws.Range('E8').Value = '  +0.43%  '
$ws.Range('D9').Value = '0.06347'
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('E10').Value = '  -0.13%  '
$ws.Range('D11').Value = '0.07830'
$ws.Range('E11').Value = '  +0.72%  '
$ws.Range('D12').Value = '4.541'
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('D13').Value = '1.657.80'
$ws.Range('E13').Value = '  +0.02%  '
$ws.Range('D14').Value = '1.887.20'
$ws.Range('E14').Value = '  +0.24%  '
$ws.Range('D15').Value = '0.5525'
$ws.Range('E15').Value = '  +0.77%  '
$ws.Range('D16').Value = '0.0₅8167'
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('D17').Value = '65.53'
$ws.Range('E17').Value = '  +0.43%  '
$ws.Range('D19').Value = '4.639'
$ws.Range('E19').Value = '  +1.81%  '
$ws.Range('D20').Value = '192.16'
$ws.Range('E20').Value = '  -0.51%  '
$ws.Range('E21').Value = '  +0.84%  '
$ws.Range('D22').Value = '6.047'
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').Value = '1.010'
$ws.Range('E23').Value = '  +0.68%  '
$ws.Range('D24').Value = '144.09'
$ws.Range('E24').Value = '  +2.82%  '
$ws.Range('E25').Value = '  -2.06%  '
$ws.Range('D26').Value = '7.218'
$ws.Range('E26').Value = '  -0.75%  '
$ws.Range('E27').Value = '  -0.34%  '
$ws.Range('D28').Value = '1.479'
$ws.Range('E28').Value = '  +3.21%  '
$ws.Range('D29').Value = '0.05857'
$ws.Range('E29').Value = '  -1.51%  '
$ws.Range('D30').Value = '1.278'
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('D31').Value = '3.577'
$ws.Range('E31').Value = '  +1.76%  '
$ws.Range('D32').Value = '3.286'
$ws.Range('E32').Value = '  +1.37%  '
$ws.Range('D33').Value = '1.612'
$ws.Range('E33').Value = '  +3.43%  '
$ws.Range('D34').Value = '0.9585'
$ws.Range('E34').Value = '  +0.75%  '
$ws.Range('D35').Value = '2.817'
$ws.Range('E35').Value = '  +1.62%  '
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('D37').Value = '0.5812'
$ws.Range('E37').Value = '  +2.66%  '
$ws.Range('E38').Value = '  -0.55%  '
$ws.Range('D39').Value = '5.895'
$ws.Range('E39').Value = '  +0.97%  '
$ws.Range('D40').Value = '0.8516'
$ws.Range('E40').Value = '  +0.75%  '
$ws.Range('E41').Value = '  +0.61%  '
$ws.Range('D42').Value = '1.041.83'
$ws.Range('E42').Value = '  +2.70%  '
$ws.Range('D43').Value = '103.81'
$ws.Range('E43').Value = '  +2.19%  '
$ws.Range('D44').Value = '1.799.88'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '0.0₈108'
$ws.Range('E45').Value = '  +4.26%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '57.07'
$ws.Range('E46').Value = '  -0.16%  '
$ws.Range('D47').Value = '1.013'
$ws.Range('E47').Value = '  +1.19%  '
$ws.Range('D48').Value = '0.4372'
$ws.Range('E48').Value = '  +2.02%  '
$ws.Range('D49').Value = '7.978'
$ws.Range('E49').Value = '  +2.47%  '
$ws.Range('E50').Value = '  +0.12%  '
$ws.Range('D51').Value = '1.428'
$ws.Range('E51').Value = '  -3.13%  '

# Restore the default (unstyled) cell style now that the text values are set,
# so no new style survives in the saved workbook.
$textRange.Style = "Normal"
